$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 183.027791
$ws.Range("D2").Value = 5.839491
$ws.Range("E2").Value = 0.016475

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 368.003
$ws.Range("D3").Value = 5.870558
$ws.Range("E3").Value = 0.003276

# Row 4 - Residuals
$ws.Range("B4").Value = 6989.512267
$ws.Range("C4").Value = 223

# Row 5 - SM-Control
$ws.Range("G5").Value = -1.729712
$ws.Range("H5").Value = -4.174776
$ws.Range("I5").Value = 0.715351
$ws.Range("J5").Value = 0.219368

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 1.106347
$ws.Range("H6").Value = -1.526469
$ws.Range("I6").Value = 3.739162
$ws.Range("J6").Value = 0.583008

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 2.836059
$ws.Range("H7").Value = 0.845623
$ws.Range("I7").Value = 4.826495
$ws.Range("J7").Value = 0.002613
